# Insert a new weekly price entry as row 20, pushing all subsequent rows
# (previous rows 20-36) down by one (to rows 21-37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 20 - shifts existing row 20..36 down to 21..37
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new data point
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C20").Value = "Los Lagos"
$ws.Range("D20").Value = 44873
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = 100112013
$ws.Range("G20").Value = "Alcachofa"
$ws.Range("H20").Value = "Española"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = 12000
$ws.Range("L20").Value = 12000
$ws.Range("M20").Value = 12000
$ws.Range("N20").Value = "$/caja 30 unidades"
$ws.Range("O20").Value = "Provincia de Limarí"
$ws.Range("P20").Value = 400
$ws.Range("Q20").Value = 30
$ws.Range("R20").Value = "Hortaliza"
